# Update specific numeric values in Sheet1 as described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D5").Value  = -7.890099999999994
$ws.Range("D6").Value  = -8.817699999999995
$ws.Range("C7").Value  = -12.06950000000001
$ws.Range("A10").Value = -20.41599999999998
$ws.Range("A12").Value = -22.27780000000002
$ws.Range("E12").Value = 12.73539999999998
$ws.Range("B13").Value = 5.926899999999998
$ws.Range("A18").Value = -22.31750000000002
$ws.Range("C20").Value = -14.58060000000001
$ws.Range("E20").Value = 11.686
$ws.Range("E25").Value = 13.05429999999999
